# Update "想去人数" (F column) figures across the four sheets to match the
# latest scrape snapshot (commit 456a3b4).

$wb = $excel.ActiveWorkbook

function Set-FValues($SheetName, $RowToValue) {
    $ws = $wb.Worksheets.Item($SheetName)
    foreach ($row in $RowToValue.Keys) {
        $ws.Cells.Item($row, 6).Value = $RowToValue[$row]
    }
}

# Sheet: 展览
Set-FValues "展览" @{
    6  = 375
    7  = 1198
    8  = 451
    9  = 7382
    12 = 2059
    13 = 8038
    14 = 46
    15 = 57
    16 = 5527
    18 = 2453
    21 = 307
    23 = 82
    24 = 13
    25 = 398
    26 = 268
    27 = 14
    28 = 2457
    30 = 274
    31 = 90
    32 = 169
    33 = 606
    36 = 1539
    38 = 8
    39 = 2406
    40 = 2227
    42 = 16
    43 = 16
}

# Sheet: 演出
Set-FValues "演出" @{
    2 = 105
    4 = 76
    5 = 18
}

# Sheet: 本地生活
Set-FValues "本地生活" @{
    3 = 1283
}

# Sheet: 全部类型
Set-FValues "全部类型" @{
    4  = 1283
    6  = 105
    7  = 375
    8  = 1198
    9  = 451
    10 = 7382
    13 = 2059
    14 = 8038
    15 = 46
    16 = 57
    17 = 5527
    19 = 2453
    23 = 82
    25 = 13
    26 = 76
    27 = 398
    28 = 268
    29 = 14
    30 = 2457
    32 = 274
    33 = 90
    34 = 169
    35 = 18
    36 = 606
    40 = 1539
    42 = 8
    43 = 2406
    45 = 2227
    47 = 16
    48 = 16
}
